$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '58.268.13'
$ws.Range("E2").Value = '  +3.29%  '
$ws.Range("D3").Value = '3.061.85'
$ws.Range("E3").Value = '  +2.70%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '522.52'
$ws.Range("E5").Value = '  +4.08%  '
$ws.Range("D6").Value = '142.56'
$ws.Range("E6").Value = '  +6.66%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +4.68%  '
$ws.Range("D9").Value = '7.51'
$ws.Range("E9").Value = '  +2.73%  '
$ws.Range("E10").Value = '  +6.29%  '
$ws.Range("E11").Value = '  +5.74%  '
$ws.Range("D12").Value = '3.568.12'
$ws.Range("E12").Value = '  +2.24%  '
$ws.Range("E13").Value = '  +2.37%  '
$ws.Range("D14").Value = '26.80'
$ws.Range("E14").Value = '  +6.82%  '
$ws.Range("D15").Value = '0.0000172'
$ws.Range("E15").Value = '  +15.09%  '
$ws.Range("D16").Value = '58.223.94'
$ws.Range("E16").Value = '  +3.23%  '
$ws.Range("D17").Value = '6.23'
$ws.Range("E17").Value = '  +9.67%  '
$ws.Range("D18").Value = '3.071.93'
$ws.Range("E18").Value = '  +3.07%  '
$ws.Range("D19").Value = '13.11'
$ws.Range("E19").Value = '  +6.24%  '
$ws.Range("E20").Value = '  +5.23%  '
$ws.Range("D21").Value = '339.06'
$ws.Range("E21").Value = '  +4.33%  '
$ws.Range("D22").Value = '5.79'
$ws.Range("E22").Value = '  +1.66%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("E24").Value = '  +7.01%  '
$ws.Range("D25").Value = '65.46'
$ws.Range("E25").Value = '  +5.30%  '
$ws.Range("E26").Value = '  +4.07%  '
$ws.Range("D27").Value = '0.0₃0966'
$ws.Range("E27").Value = '  +8.74%  '
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").Value = '6.96'
$ws.Range("E29").Value = '  +7.73%  '
$ws.Range("D30").Value = '7.54'
$ws.Range("E30").Value = '  +11.34%  '
$ws.Range("E31").Value = '  +5.79%  '
$ws.Range("D32").Value = '1.23'
$ws.Range("E32").Value = '  +4.73%  '
$ws.Range("D33").Value = '21.19'
$ws.Range("E33").Value = '  +4.75%  '
$ws.Range("E34").Value = '  +8.14%  '
$ws.Range("D35").Value = '157.62'
$ws.Range("E35").Value = '  +0.95%  '
$ws.Range("D36").Value = '5.96'
$ws.Range("E36").Value = '  +8.01%  '
$ws.Range("E37").Value = '  +2.59%  '
$ws.Range("D38").Value = '25.48'
$ws.Range("E38").Value = '  +10.69%  '
$ws.Range("D39").Value = '0.0695'
$ws.Range("E39").Value = '  +3.53%  '
$ws.Range("D40").Value = '3.100.10'
$ws.Range("E40").Value = '  +2.85%  '
$ws.Range("D41").Value = '37.74'
$ws.Range("E41").Value = '  +4.88%  '
$ws.Range("E42").Value = '  +10.12%  '
$ws.Range("E43").Value = '  +4.13%  '
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").Value = '2.335.80'
$ws.Range("E45").Value = '  +4.48%  '
$ws.Range("E46").Value = '  +4.85%  '
$ws.Range("E47").Value = '  +3.42%  '
$ws.Range("D48").Value = '6.09'
$ws.Range("E48").Value = '  +5.85%  '
$ws.Range("E49").Value = '  +3.16%  '
$ws.Range("D50").Value = '19.88'
$ws.Range("E50").Value = '  +5.24%  '
$ws.Range("E51").Value = '  -1.72%  '
